$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-07 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-08 Friday", 2) | Out-Null
$d.Content.Find.Execute("123÷9=13, 6", $true, $false, $false, $false, $false, $true, 1, $false, "361÷8=45, 1", 2) | Out-Null
$d.Content.Find.Execute("508÷7=72, 4", $true, $false, $false, $false, $false, $true, 1, $false, "252÷7=36, 0", 2) | Out-Null
$d.Content.Find.Execute("159÷3=53, 0", $true, $false, $false, $false, $false, $true, 1, $false, "845÷2=422, 1", 2) | Out-Null
$d.Content.Find.Execute("329÷8=41, 1", $true, $false, $false, $false, $false, $true, 1, $false, "975÷7=139, 2", 2) | Out-Null
$d.Content.Find.Execute("638÷2=319, 0", $true, $false, $false, $false, $false, $true, 1, $false, "111÷4=27, 3", 2) | Out-Null
$d.Content.Find.Execute("887÷9=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "552÷9=61, 3", 2) | Out-Null
$d.Content.Find.Execute("425÷8=53, 1", $true, $false, $false, $false, $false, $true, 1, $false, "381÷7=54, 3", 2) | Out-Null
$d.Content.Find.Execute("853÷6=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "828÷8=103, 4", 2) | Out-Null
$d.Content.Find.Execute("881÷9=97, 8", $true, $false, $false, $false, $false, $true, 1, $false, "837÷5=167, 2", 2) | Out-Null
$d.Content.Find.Execute("597÷9=66, 3", $true, $false, $false, $false, $false, $true, 1, $false, "214÷5=42, 4", 2) | Out-Null
$d.Content.Find.Execute("686÷8=85, 6", $true, $false, $false, $false, $false, $true, 1, $false, "973÷3=324, 1", 2) | Out-Null
$d.Content.Find.Execute("731÷6=121, 5", $true, $false, $false, $false, $false, $true, 1, $false, "307÷6=51, 1", 2) | Out-Null
$d.Content.Find.Execute("385÷6=64, 1", $true, $false, $false, $false, $false, $true, 1, $false, "675÷8=84, 3", 2) | Out-Null
$d.Content.Find.Execute("794÷8=99, 2", $true, $false, $false, $false, $false, $true, 1, $false, "414÷9=46, 0", 2) | Out-Null
$d.Content.Find.Execute("515÷9=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "801÷3=267, 0", 2) | Out-Null
$d.Content.Find.Execute("142÷4=35, 2", $true, $false, $false, $false, $false, $true, 1, $false, "288÷3=96, 0", 2) | Out-Null
$d.Content.Find.Execute("342÷4=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "858÷8=107, 2", 2) | Out-Null
$d.Content.Find.Execute("747÷7=106, 5", $true, $false, $false, $false, $false, $true, 1, $false, "640÷8=80, 0", 2) | Out-Null
$d.Content.Find.Execute("882÷3=294, 0", $true, $false, $false, $false, $false, $true, 1, $false, "265÷4=66, 1", 2) | Out-Null
$d.Content.Find.Execute("673÷9=74, 7", $true, $false, $false, $false, $false, $true, 1, $false, "985÷8=123, 1", 2) | Out-Null
$d.Content.Find.Execute("468÷4=117, 0", $true, $false, $false, $false, $false, $true, 1, $false, "534÷9=59, 3", 2) | Out-Null
$d.Content.Find.Execute("656÷3=218, 2", $true, $false, $false, $false, $false, $true, 1, $false, "965÷9=107, 2", 2) | Out-Null
$d.Content.Find.Execute("454÷7=64, 6", $true, $false, $false, $false, $false, $true, 1, $false, "124÷2=62, 0", 2) | Out-Null
$d.Content.Find.Execute("263÷8=32, 7", $true, $false, $false, $false, $false, $true, 1, $false, "214÷8=26, 6", 2) | Out-Null
$d.Content.Find.Execute("432÷9=48, 0", $true, $false, $false, $false, $false, $true, 1, $false, "777÷3=259, 0", 2) | Out-Null
